# The workbook opens with "Editorial_Rank" already the active sheet/tab.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data edit: B2 ("wallpaper.com" editorial rank) changes from 4 to 3.
$ws.Range("B2").Value = 3

# Move/record the cell selection onto B2 (previously C9), matching the
# saved sheetView's <selection activeCell="B2" sqref="B2"/>.
[void]$ws.Range("B2").Select()
